# The import file now only contains rows up to r=385 (rows 386-402 were
# removed, shrinking the used range from A1:J402 to A1:J385).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the trailing rows (386 through 402) that are no longer part of
# the imported data.
$ws.Range("A386:J402").EntireRow.Delete()
